# Performance_KPIs.xlsx - "First batch of bug fixes"
# The underlying financial model was recalculated with corrected inputs; this
# script writes the refreshed literal results into both report sheets and
# drops the now-removed "year 30" row from the CoCRoI sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CoCRoI")
$ws2 = $wb.Worksheets.Item("Overall CAGR")

# --- CoCRoI sheet: refreshed "CoCRoI" (B) and "Annual Cash Flow" (C) values for years 0-29 ---
$ws1.Range("B2").Value = -9.03663356112053
$ws1.Range("C2").Value = -3963.693395746492
$ws1.Range("B3").Value = -3.680311722040075
$ws1.Range("C3").Value = -1614.276729079828
$ws1.Range("B4").Value = -3.288989407990486
$ws1.Range("C4").Value = -1442.632979079827
$ws1.Range("B5").Value = -2.888595349854266
$ws1.Range("C5").Value = -1267.010135329828
$ws1.Range("B6").Value = -2.478931206579836
$ws1.Range("C6").Value = -1087.321200486081
$ws1.Range("B7").Value = -2.059794816691333
$ws1.Range("C7").Value = -903.4775014712359
$ws1.Range("B8").Value = -1.630980148302129
$ws1.Range("C8").Value = -715.3886675490213
$ws1.Range("B9").Value = -1.192277249699567
$ws1.Range("C9").Value = -522.9626086494727
$ws1.Range("B10").Value = -0.743472200588171
$ws1.Range("C10").Value = -326.1054939829865
$ws1.Range("B11").Value = -0.2843470640834478
$ws1.Range("C11").Value = -124.7217309836023
$ws1.Range("B12").Value = 0.1853201604460281
$ws1.Range("C12").Value = 81.28605537563908
$ws1.Range("B13").Value = 0.6657555835853834
$ws1.Range("C13").Value = 292.0170428501388
$ws1.Range("B14").Value = 1.157189471019717
$ws1.Range("C14").Value = 507.5722317260233
$ws1.Range("B15").Value = 1.659856287505067
$ws1.Range("C15").Value = 728.05446410691
$ws1.Range("B16").Value = 2.173994739542318
$ws1.Range("C16").Value = 953.5684426317493
$ws1.Range("B17").Value = 2.699847816625867
$ws1.Range("C17").Value = 1184.220748567521
$ws1.Range("B18").Value = 3.237662830931656
$ws1.Range("C18").Value = 1420.119859217398
$ws1.Range("B19").Value = 3.787691455302082
$ws1.Range("C19").Value = 1661.376164581875
$ws1.Range("B20").Value = 4.350189759377022
$ws1.Range("C20").Value = 1908.101983206746
$ws1.Range("B21").Value = 4.92541824371289
$ws1.Range("C21").Value = 2160.411577148567
$ws1.Range("B22").Value = 5.513641871722593
$ws1.Range("C22").Value = 2418.421165984323
$ws1.Range("B23").Value = 6.115130099260587
$ws1.Range("C23").Value = 2682.248939788175
$ws1.Range("B24").Value = 6.730156901668202
$ws1.Range("C24").Value = 2952.015070994215
$ws1.Range("B25").Value = 7.359000798084413
$ws1.Range("C25").Value = 3227.841725059776
$ws1.Range("B26").Value = 8.001944872817386
$ws1.Range("C26").Value = 3509.853069839526
$ws1.Range("B27").Value = 8.659276793561284
$ws1.Range("C27").Value = 3798.175283575819
$ws1.Range("B28").Value = 9.331288826232104
$ws1.Range("C28").Value = 4092.936561406057
$ws1.Range("B29").Value = 10.01827784618431
$ws1.Range("C29").Value = 4394.267120282594
$ws1.Range("B30").Value = 10.72054534555858
$ws1.Range("C30").Value = 4702.299202195632
$ws1.Range("B31").Value = 11.43839743649735
$ws1.Range("C31").Value = 5017.167075583651

# Year 30 (row 32) no longer exists in the refreshed model -> drop the row
# (this also shrinks the sheet dimension from A1:D32 to A1:D31).
$ws1.Rows.Item(32).Delete()

# --- Overall CAGR sheet: refreshed "Cumulative Annual Cash Flow" (B), "Annual Cash Flow" (C),
#     "Overall Return" (G) and "Overall CAGR" (H) values for years 0-29 ---
$ws2.Range("B2").Value = -3963.693395746492
$ws2.Range("C2").Value = -3963.693395746492
$ws2.Range("G2").Value = 21764.6591353917
$ws2.Range("B3").Value = -5577.97012482632
$ws2.Range("C3").Value = -1614.276729079828
$ws2.Range("G3").Value = 23976.12694880568
$ws2.Range("H3").Value = -45.33798358778984
$ws2.Range("B4").Value = -7020.603103906147
$ws2.Range("C4").Value = -1442.632979079827
$ws2.Range("G4").Value = 26497.83563571166
$ws2.Range("H4").Value = -22.27539737932498
$ws2.Range("B5").Value = -8287.613239235974
$ws2.Range("C5").Value = -1267.010135329828
$ws2.Range("G5").Value = 29341.10693919603
$ws2.Range("H5").Value = -12.54304124930538
$ws2.Range("B6").Value = -9374.934439722056
$ws2.Range("C6").Value = -1087.321200486081
$ws2.Range("G6").Value = 32517.83296379016
$ws2.Range("H6").Value = -7.208743916184678
$ws2.Range("B7").Value = -10278.41194119329
$ws2.Range("C7").Value = -903.4775014712359
$ws2.Range("G7").Value = 36040.51229358723
$ws2.Range("H7").Value = -3.85216381063862
$ws2.Range("B8").Value = -10993.80060874231
$ws2.Range("C8").Value = -715.3886675490213
$ws2.Range("G8").Value = 39922.2886451615
$ws2.Range("H8").Value = -1.556508558850656
$ws2.Range("B9").Value = -11516.76321739178
$ws2.Range("C9").Value = -522.9626086494727
$ws2.Range("G9").Value = 44176.99223958647
$ws2.Range("H9").Value = 0.1021146026615849
$ws2.Range("B10").Value = -11842.86871137477
$ws2.Range("C10").Value = -326.1054939829865
$ws2.Range("G10").Value = 48819.18409138946
$ws2.Range("H10").Value = 1.347290393893985
$ws2.Range("B11").Value = -11967.59044235837
$ws2.Range("C11").Value = -124.7217309836023
$ws2.Range("G11").Value = 53864.20342681867
$ws2.Range("H11").Value = 2.308536771441183
$ws2.Range("B12").Value = -11886.30438698273
$ws2.Range("C12").Value = 81.28605537563908
$ws2.Range("G12").Value = 59328.21845940418
$ws2.Range("H12").Value = 3.066325420585914
$ws2.Range("B13").Value = -11594.2873441326
$ws2.Range("C13").Value = 292.0170428501388
$ws2.Range("G13").Value = 65228.28076755338
$ws2.Range("H13").Value = 3.673438915967031
$ws2.Range("B14").Value = -11086.71511240657
$ws2.Range("C14").Value = 507.5722317260233
$ws2.Range("G14").Value = 71582.38353691186
$ws2.Range("H14").Value = 4.166018479049161
$ws2.Range("B15").Value = -10358.66064829966
$ws2.Range("C15").Value = 728.05446410691
$ws2.Range("G15").Value = 78409.52394953823
$ws2.Range("H15").Value = 4.569685718763727
$ws2.Range("B16").Value = -9405.092205667912
$ws2.Range("C16").Value = 953.5684426317493
$ws2.Range("G16").Value = 85729.77002268052
$ws2.Range("H16").Value = 4.903130468032946
$ws2.Range("B17").Value = -8220.871457100391
$ws2.Range("C17").Value = 1184.220748567521
$ws2.Range("G17").Value = 93564.33222220994
$ws2.Range("H17").Value = 5.180313937185099
$ws2.Range("B18").Value = -6800.751597882993
$ws2.Range("C18").Value = 1420.119859217398
$ws2.Range("G18").Value = 101935.6401996741
$ws2.Range("H18").Value = 5.411875908373243
$ws2.Range("B19").Value = -5139.375433301118
$ws2.Range("C19").Value = 1661.376164581875
$ws2.Range("G19").Value = 110867.4250276026
$ws2.Range("H19").Value = 5.606064270430533
$ws2.Range("B20").Value = -3231.273450094372
$ws2.Range("C20").Value = 1908.101983206746
$ws2.Range("G20").Value = 120384.8073352547
$ws2.Range("H20").Value = 5.769367090352961
$ws2.Range("B21").Value = -1070.861872945805
$ws2.Range("C21").Value = 2160.411577148567
$ws2.Range("G21").Value = 130514.3917765912
$ws2.Range("H21").Value = 5.906953376021495
$ws2.Range("B22").Value = 1347.559293038517
$ws2.Range("C22").Value = 2418.421165984323
$ws2.Range("G22").Value = 141284.36829402
$ws2.Range("H22").Value = 6.022987279868697
$ws2.Range("B23").Value = 4029.808232826692
$ws2.Range("C23").Value = 2682.248939788175
$ws2.Range("G23").Value = 152724.6206755754
$ws2.Range("H23").Value = 6.120856467677505
$ws2.Range("B24").Value = 6981.823303820907
$ws2.Range("C24").Value = 2952.015070994215
$ws2.Range("G24").Value = 164866.842939813
$ws2.Range("H24").Value = 6.203340972578109
$ws2.Range("B25").Value = 10209.66502888068
$ws2.Range("C25").Value = 3227.841725059776
$ws2.Range("G25").Value = 177744.6641220223
$ws2.Range("H25").Value = 6.272739963310681
$ws2.Range("B26").Value = 13719.51809872021
$ws2.Range("C26").Value = 3509.853069839526
$ws2.Range("G26").Value = 191393.7820775753
$ws2.Range("H26").Value = 6.330968222576416
$ws2.Range("B27").Value = 17517.69338229603
$ws2.Range("C27").Value = 3798.175283575819
$ws2.Range("G27").Value = 205852.1069635561
$ws2.Range("H27").Value = 6.379630477089204
$ws2.Range("B28").Value = 21610.62994370209
$ws2.Range("C28").Value = 4092.936561406057
$ws2.Range("G28").Value = 221159.9151084815
$ws2.Range("H28").Value = 6.420079299425785
$ws2.Range("B29").Value = 26004.89706398468
$ws2.Range("C29").Value = 4394.267120282594
$ws2.Range("G29").Value = 237360.0140321722
$ws2.Range("H29").Value = 6.453460665729449
$ws2.Range("B30").Value = 30707.19626618031
$ws2.Range("C30").Value = 4702.299202195632
$ws2.Range("G30").Value = 254497.9194339343
$ws2.Range("H30").Value = 6.480750128269364
$ws2.Range("B31").Value = 35724.36334176396
$ws2.Range("C31").Value = 5017.167075583651
$ws2.Range("G31").Value = 272622.0450274411
$ws2.Range("H31").Value = 6.502781775600242

